# weight_over_time.xlsx - add last weekend's (disappointing, first keto week)
# weigh-in measurements to the raw_data sheet (rows 82-90).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# New rows of date/time/weight data, most-recent-first (matching existing order)
# columns: row, dateSerial (A, date+time), timeFraction (B, time of day), weight (C)
$newRows = @(
    @{ Row = 82; Date = 44081.324305555558; Time = 0.32430555555555557; Weight = 73.099999999999994 },
    @{ Row = 83; Date = 44081.314583333333; Time = 0.31458333333333333; Weight = 73.099999999999994 },
    @{ Row = 84; Date = 44080.918055555558; Time = 0.91805555555555562; Weight = 74.5 },
    @{ Row = 85; Date = 44080.344444444447; Time = 0.3444444444444445;  Weight = 74 },
    @{ Row = 86; Date = 44080.341666666667; Time = 0.34166666666666662; Weight = 74.5 },
    @{ Row = 87; Date = 44079.913888888892; Time = 0.91388888888888886; Weight = 74.5 },
    @{ Row = 88; Date = 44079.297222222223; Time = 0.29722222222222222; Weight = 74.5 },
    @{ Row = 89; Date = 44078.893055555556; Time = 0.8930555555555556;  Weight = 74.5 },
    @{ Row = 90; Date = 44078.342361111114; Time = 0.34236111111111112; Weight = 72.599999999999994 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    # Column A: date+time, formatted like the rest of the recent rows (m/d/yyyy h:mm)
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 1).NumberFormat = "m/d/yy h:mm"

    # Column B: time-of-day fraction, formatted as h:mm
    $ws.Cells.Item($row, 2).Value = $r.Time
    $ws.Cells.Item($row, 2).NumberFormat = "h:mm"

    # Column C: weight (plain number)
    $ws.Cells.Item($row, 3).Value = $r.Weight
}

# Row 87 in the source workbook carries slightly different (copy/pasted) styling:
# the date cell uses the short-date format and the weight cell inherited a time format.
$ws.Cells.Item(87, 1).NumberFormat = "m/d/yy"
$ws.Cells.Item(87, 3).NumberFormat = "h:mm"

# Extend the shared "AM"/"PM" formula down through the new rows (D41:D90)
$ws.Range("D41:D90").Formula = '=IF(B41<TIME(12,0,0), "AM", "PM")'

# Update the active selection to reflect the new end of the data (A91, one past last row)
$ws.Activate()
$ws.Range("A91").Select()
